$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "kahal"
$ws.Range("C2").Value = "kahal@gmail.com"

# Update row 3
$ws.Range("A3").Value = 6
$ws.Range("B3").Value = "kahal"
$ws.Range("C3").Value = "kaaaaaaaahal@gmail.com"

# Add row 4
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "omarrrrrrrrrrrrrrr"
$ws.Range("C4").Value = "wq@gmail.com"
